$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1228.6923
$ws.Range("I12").Value = 1849.5
$ws.Range("J12").Value = 952.7778
$ws.Range("K12").Value = 1849.5
$ws.Range("L12").Value = 952.7778
$ws.Range("M12").Value = -1679.5
$ws.Range("N12").Value = -1292.7778
$ws.Range("H51").Value = 7499.5
$ws.Range("I51").Value = 7499.5
$ws.Range("K51").Value = 7499.5
$ws.Range("M51").Value = -7015.5
$ws.Range("H92").Value = 584
$ws.Range("I92").Value = 584
$ws.Range("K92").Value = 584
$ws.Range("M92").Value = 664
$ws.Range("H97").Value = 250001380
$ws.Range("J97").Value = 250001380
$ws.Range("L97").Value = 750004140
$ws.Range("N97").Value = -750005132
$ws.Range("H99").Value = 1504.6666
$ws.Range("I99").Value = 1504.6666
$ws.Range("K99").Value = 4513.9998
$ws.Range("M99").Value = -3015.9998
$ws.Range("H116").Value = 7439.2
$ws.Range("I116").Value = 6698.25
$ws.Range("J116").Value = 7933.1665
$ws.Range("K116").Value = 6698.25
$ws.Range("L116").Value = 7933.1665
$ws.Range("M116").Value = -3256.25
$ws.Range("N116").Value = -14817.1665
$ws.Range("H135").Value = 945.5
$ws.Range("I135").Value = 945.5
$ws.Range("K135").Value = 8509.5
$ws.Range("M135").Value = -5974.5
$ws.Range("H137").Value = 2852.5715
$ws.Range("J137").Value = 3197.8
$ws.Range("L137").Value = 9593.400000000001
$ws.Range("N137").Value = -14693.4
$ws.Range("H138").Value = 14842.053
$ws.Range("I138").Value = 9000
$ws.Range("J138").Value = 15166.611
$ws.Range("K138").Value = 27000
$ws.Range("L138").Value = 45499.833
$ws.Range("M138").Value = -21860
$ws.Range("N138").Value = -55779.833

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1074.5
$ws.Range("I5").Value = 1999
$ws.Range("K5").Value = 1999
$ws.Range("M5").Value = -1887
$ws.Range("H11").Value = 1099.5
$ws.Range("J11").Value = 1099.5
$ws.Range("L11").Value = 1099.5
$ws.Range("N11").Value = -1387.5
$ws.Range("H45").Value = 2467.1428
$ws.Range("I45").Value = 1726.6666
$ws.Range("K45").Value = 1726.6666
$ws.Range("M45").Value = -1349.6666
$ws.Range("H74").Value = 1984.8572
$ws.Range("I74").Value = 1984.8572
$ws.Range("K74").Value = 1984.8572
$ws.Range("M74").Value = -1110.8572
$ws.Range("H77").Value = 1984.8572
$ws.Range("I77").Value = 1984.8572
$ws.Range("K77").Value = 9924.286
$ws.Range("M77").Value = -5556.286
$ws.Range("H98").Value = 27298
$ws.Range("J98").Value = 27298
$ws.Range("L98").Value = 27298
$ws.Range("N98").Value = -33288
$ws.Range("H114").Value = 75000
$ws.Range("J114").Value = 75000
$ws.Range("L114").Value = 75000
$ws.Range("N114").Value = -83678

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1074.5
$ws.Range("I4").Value = 1999
$ws.Range("K4").Value = 1999
$ws.Range("M4").Value = -1884
$ws.Range("H22").Value = 831
$ws.Range("I22").Value = 831
$ws.Range("K22").Value = 831
$ws.Range("M22").Value = -658
$ws.Range("H94").Value = 3163.7188
$ws.Range("I94").Value = 1725.7368
$ws.Range("J94").Value = 5265.385
$ws.Range("K94").Value = 1725.7368
$ws.Range("L94").Value = 5265.385
$ws.Range("M94").Value = -1274.7368
$ws.Range("N94").Value = -6167.385

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6736.875
$ws.Range("I99").Value = 4428.7144
$ws.Range("K99").Value = 4428.7144
$ws.Range("M99").Value = -2930.7144
$ws.Range("H126").Value = 6736.875
$ws.Range("I126").Value = 4428.7144
$ws.Range("K126").Value = 13286.1432
$ws.Range("M126").Value = -10816.1432
$ws.Range("H132").Value = 3245.3157
$ws.Range("J132").Value = 15000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -50060
$ws.Range("H134").Value = 5329.2383
$ws.Range("I134").Value = 5212
$ws.Range("K134").Value = 15636
$ws.Range("M134").Value = -13101
$ws.Range("H141").Value = 470446.72
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 570546
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 570546
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -580906

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 187.875
$ws.Range("I6").Value = 216.83333
$ws.Range("K6").Value = 650.49999
$ws.Range("M6").Value = -537.49999
$ws.Range("H15").Value = 280.5
$ws.Range("I15").Value = 268.66666
$ws.Range("J15").Value = 298.25
$ws.Range("K15").Value = 805.9999799999999
$ws.Range("L15").Value = 894.75
$ws.Range("M15").Value = -665.9999799999999
$ws.Range("N15").Value = -1174.75
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H55").Value = 2808.125
$ws.Range("J55").Value = 2749.1667
$ws.Range("L55").Value = 8247.500100000001
$ws.Range("N55").Value = -8601.500100000001
$ws.Range("H107").Value = 1911.2778
$ws.Range("J107").Value = 2212.6667
$ws.Range("L107").Value = 6638.000100000001
$ws.Range("N107").Value = -10478.0001
$ws.Range("H132").Value = 1798.4
$ws.Range("I132").Value = 1998
$ws.Range("J132").Value = 1748.5
$ws.Range("K132").Value = 17982
$ws.Range("L132").Value = 15736.5
$ws.Range("M132").Value = -15452
$ws.Range("N132").Value = -20796.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1150.125
$ws.Range("I97").Value = 743.2857
$ws.Range("K97").Value = 743.2857
$ws.Range("M97").Value = -247.2857
$ws.Range("H102").Value = 1974.75
$ws.Range("I102").Value = 950
$ws.Range("K102").Value = 950
$ws.Range("M102").Value = 672
$ws.Range("H122").Value = 4185.9287
$ws.Range("I122").Value = 4412.231
$ws.Range("K122").Value = 13236.693
$ws.Range("M122").Value = -10786.693
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 112348490
$ws.Range("I126").Value = 280865470
$ws.Range("J126").Value = 3841.1667
$ws.Range("K126").Value = 842596410
$ws.Range("L126").Value = 11523.5001
$ws.Range("M126").Value = -842593940
$ws.Range("N126").Value = -16463.5001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3136.125
$ws.Range("I7").Value = 3155.5715
$ws.Range("K7").Value = 3155.5715
$ws.Range("M7").Value = -3043.5715
$ws.Range("H22").Value = 644
$ws.Range("I22").Value = 593.5
$ws.Range("K22").Value = 593.5
$ws.Range("M22").Value = -298.5
$ws.Range("H27").Value = 644
$ws.Range("I27").Value = 593.5
$ws.Range("K27").Value = 593.5
$ws.Range("M27").Value = -486.5
$ws.Range("H40").Value = 50000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 50000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 50000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -50272
$ws.Range("H103").Value = 14999.5
$ws.Range("J103").Value = 14999.5
$ws.Range("L103").Value = 14999.5
$ws.Range("N103").Value = -17343.5
$ws.Range("H122").Value = 7493.7
$ws.Range("I122").Value = 7493.7
$ws.Range("K122").Value = 22481.1
$ws.Range("M122").Value = -20031.1
$ws.Range("H126").Value = 3136.125
$ws.Range("I126").Value = 3155.5715
$ws.Range("K126").Value = 9466.7145
$ws.Range("M126").Value = -6996.7145
$ws.Range("H132").Value = 9238.799999999999
$ws.Range("I132").Value = 6770.4287
$ws.Range("K132").Value = 20311.2861
$ws.Range("M132").Value = -17781.2861
$ws.Range("H136").Value = 1949.5
$ws.Range("I136").Value = 1949.5
$ws.Range("K136").Value = 5848.5
$ws.Range("M136").Value = -3298.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10002
$ws.Range("I62").Value = 10001
$ws.Range("K62").Value = 10001
$ws.Range("M62").Value = -9377
$ws.Range("H65").Value = 10002
$ws.Range("I65").Value = 10001
$ws.Range("K65").Value = 50005
$ws.Range("M65").Value = -46885
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H126").Value = 1386.091
$ws.Range("I126").Value = 1414.7
$ws.Range("K126").Value = 4244.1
$ws.Range("M126").Value = -1774.1
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 44265.918
$ws.Range("I136").Value = 51619.1
$ws.Range("K136").Value = 154857.3
$ws.Range("M136").Value = -152307.3
